$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("G4").Value = 32
$ws.Range("J4").Value = 23
$ws.Range("L4").Value = 1.23
$ws.Range("Q4").Value = 1.16
$ws.Range("R4").Value = 3.88
$ws.Range("U4").Value = 2.88
$ws.Range("V4").Value = 1.39
$ws.Range("W4").Value = 150
$ws.Range("Y4").Value = 175
$ws.Range("AB4").Value = 500
$ws.Range("AC4").Value = 25
$ws.Range("AD4").Value = 29
$ws.Range("AF4").Value = 300
$ws.Range("AH4").Value = 12
$ws.Range("AI4").Value = 6.8
$ws.Range("AK4").Value = 5.4
$ws.Range("AL4").Value = 12
$ws.Range("AM4").Value = 50
$ws.Range("AN4").Value = 35
$ws.Range("AO4").Value = 400
$ws.Range("AP4").Value = 175
$ws.Range("AT4").Value = 5.2
$ws.Range("AU4").Value = 15.5
$ws.Range("AW4").Value = 3.25
$ws.Range("AX4").Value = 3.65
$ws.Range("AY4").Value = 17.5
$ws.Range("AZ4").Value = 6
$ws.Range("BA4").Value = 28
$ws.Range("BB4").Value = 250

# Row 5 updates
$ws.Range("O5").Value = 1.15
$ws.Range("P5").Value = 3.98

# Row 7 updates
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("O7").Value = 1.33
$ws.Range("P7").Value = 3.25
$ws.Range("Q7").Value = 2.08
$ws.Range("R7").Value = 1.73
